# Applies the "сумм"/"тк" column addition (I, J) to the gradebook sheet:
#  - I3 header "сумм", J3 header "тк"
#  - I4:I31 = SUM(C:H) per row
#  - J4:J31 = literal attendance-count values
#  - color-scale conditional formatting on I4:I31
#  - selection moved to J20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 3): new "сумм" / "тк" column labels
$ws.Range("I3").Value = "сумм"
$ws.Range("J3").Value = "тк"

# Literal values for column J (row 4..31), taken in order
$jValues = @(3,5,5,4,5,5,5,4,4,5,5,5,5,5,4,3,4,4,5,5,5,5,5,5,5,5,5,3)

for ($row = 4; $row -le 31; $row++) {
    $formula = "=SUM(C" + $row + ":H" + $row + ")"
    $ws.Range("I" + $row).Formula = $formula

    $jValue = $jValues[$row - 4]
    $ws.Range("J" + $row).Value = $jValue
}

# Color-scale conditional formatting over the new sum column
$cfRange = $ws.Range("I4:I31")
$cf = $cfRange.FormatConditions.AddColorScale(3)
$cf.ColorScaleCriteria.Item(1).Type = 1
$cf.ColorScaleCriteria.Item(1).FormatColor.Color = 7039083
$cf.ColorScaleCriteria.Item(2).Type = 4
$cf.ColorScaleCriteria.Item(2).Value = 50
$cf.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167
$cf.ColorScaleCriteria.Item(3).Type = 2
$cf.ColorScaleCriteria.Item(3).FormatColor.Color = 6531452

# Move the active selection to J20, matching the saved view state
[void]$ws.Range("J20").Select()
